$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1707.4166
$ws.Range("I2").Value = 1125
$ws.Range("J2").Value = 1823.9
$ws.Range("K2").Value = 1125
$ws.Range("L2").Value = 1823.9
$ws.Range("M2").Value = -1012
$ws.Range("N2").Value = -2049.9
$ws.Range("H33").Value = 249.89473
$ws.Range("I33").Value = 229.76923
$ws.Range("J33").Value = 293.5
$ws.Range("K33").Value = 229.76923
$ws.Range("L33").Value = 293.5
$ws.Range("M33").Value = -0.7692299999999932
$ws.Range("N33").Value = -751.5
$ws.Range("H64").Value = 5750
$ws.Range("I64").Value = 3500
$ws.Range("K64").Value = 3500
$ws.Range("M64").Value = -3252
$ws.Range("H67").Value = 5750
$ws.Range("I67").Value = 3500
$ws.Range("K67").Value = 3500
$ws.Range("M67").Value = -2642
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("H107").Value = 1213.2307
$ws.Range("I107").Value = 875.86957
$ws.Range("J107").Value = 3799.6667
$ws.Range("K107").Value = 875.86957
$ws.Range("L107").Value = 3799.6667
$ws.Range("M107").Value = 1044.13043
$ws.Range("N107").Value = -7639.6667
$ws.Range("H138").Value = 4188.25
$ws.Range("I138").Value = 2160.8333
$ws.Range("J138").Value = 5057.143
$ws.Range("K138").Value = 6482.499899999999
$ws.Range("L138").Value = 15171.429
$ws.Range("M138").Value = -1342.499899999999
$ws.Range("N138").Value = -25451.429
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1779.75
$ws.Range("I45").Value = 2006
$ws.Range("J45").Value = 1101
$ws.Range("K45").Value = 2006
$ws.Range("L45").Value = 1101
$ws.Range("M45").Value = -1629
$ws.Range("N45").Value = -1855
$ws.Range("H61").Value = 4266.5557
$ws.Range("I61").Value = 2715.5
$ws.Range("J61").Value = 5507.4
$ws.Range("K61").Value = 2715.5
$ws.Range("L61").Value = 5507.4
$ws.Range("M61").Value = -2503.5
$ws.Range("N61").Value = -5931.4
$ws.Range("H88").Value = 1277
$ws.Range("J88").Value = 1001
$ws.Range("L88").Value = 1001
$ws.Range("N88").Value = -1813
$ws.Range("H91").Value = 1277
$ws.Range("J91").Value = 1001
$ws.Range("L91").Value = 1001
$ws.Range("N91").Value = -3809
$ws.Range("H132").Value = 3921
$ws.Range("I132").Value = 1673.5
$ws.Range("K132").Value = 5020.5
$ws.Range("M132").Value = -2490.5
$ws.Range("H136").Value = 4266.5557
$ws.Range("I136").Value = 2715.5
$ws.Range("J136").Value = 5507.4
$ws.Range("K136").Value = 8146.5
$ws.Range("L136").Value = 16522.2
$ws.Range("M136").Value = -5596.5
$ws.Range("N136").Value = -21622.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620
$ws.Range("H80").Value = 314.25
$ws.Range("J80").Value = 415.66666
$ws.Range("L80").Value = 415.66666
$ws.Range("N80").Value = -2411.66666
$ws.Range("H83").Value = 314.25
$ws.Range("J83").Value = 415.66666
$ws.Range("L83").Value = 2078.3333
$ws.Range("N83").Value = -12062.3333
$ws.Range("H86").Value = 2445.2222
$ws.Range("J86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2445.2222
$ws.Range("J89").Value = 2000
$ws.Range("L89").Value = 10000
$ws.Range("N89").Value = -21232
$ws.Range("H134").Value = 6250.1665
$ws.Range("I134").Value = 6870.6665
$ws.Range("K134").Value = 20611.9995
$ws.Range("M134").Value = -18076.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 30090.5
$ws.Range("I82").Value = 30000
$ws.Range("K82").Value = 30000
$ws.Range("M82").Value = -29639
$ws.Range("H85").Value = 30090.5
$ws.Range("I85").Value = 30000
$ws.Range("K85").Value = 30000
$ws.Range("M85").Value = -28752
$ws.Range("H93").Value = 11400
$ws.Range("I93").Value = 8000
$ws.Range("J93").Value = 25000
$ws.Range("K93").Value = 8000
$ws.Range("L93").Value = 25000
$ws.Range("M93").Value = -6128
$ws.Range("N93").Value = -28744
$ws.Range("H99").Value = 5459.8
$ws.Range("I99").Value = 5449.75
$ws.Range("J99").Value = 5500
$ws.Range("K99").Value = 5449.75
$ws.Range("L99").Value = 5500
$ws.Range("M99").Value = -3951.75
$ws.Range("N99").Value = -8496
$ws.Range("H122").Value = 3985.762
$ws.Range("J122").Value = 4087.6667
$ws.Range("L122").Value = 12263.0001
$ws.Range("N122").Value = -17163.0001
$ws.Range("H126").Value = 5459.8
$ws.Range("I126").Value = 5449.75
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 16349.25
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -13879.25
$ws.Range("N126").Value = -21440
$ws.Range("H132").Value = 1598.2222
$ws.Range("I132").Value = 1598.2222
$ws.Range("K132").Value = 4794.6666
$ws.Range("M132").Value = -2264.6666
$ws.Range("H134").Value = 8200
$ws.Range("I134").Value = 6500
$ws.Range("K134").Value = 19500
$ws.Range("M134").Value = -16965

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1121.8889
$ws.Range("I23").Value = 1033
$ws.Range("J23").Value = 1166.3334
$ws.Range("K23").Value = 3099
$ws.Range("L23").Value = 3499.0002
$ws.Range("M23").Value = -2864
$ws.Range("N23").Value = -3969.0002
$ws.Range("H34").Value = 1581.15
$ws.Range("I34").Value = 174.81818
$ws.Range("K34").Value = 524.4545400000001
$ws.Range("M34").Value = -440.4545400000001
$ws.Range("H55").Value = 143.6
$ws.Range("J55").Value = 85
$ws.Range("L55").Value = 255
$ws.Range("N55").Value = -609
$ws.Range("H102").Value = 5000
$ws.Range("I102").Value = 5000
$ws.Range("K102").Value = 15000
$ws.Range("M102").Value = -12566

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7618.8335
$ws.Range("J70").Value = 7672.6665
$ws.Range("L70").Value = 7672.6665
$ws.Range("N70").Value = -8212.666499999999
$ws.Range("H73").Value = 7618.8335
$ws.Range("J73").Value = 7672.6665
$ws.Range("L73").Value = 7672.6665
$ws.Range("N73").Value = -9544.666499999999
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H132").Value = 1895.1666
$ws.Range("I132").Value = 1974.2
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 5922.6
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -3392.6
$ws.Range("N132").Value = -9560
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4074.2
$ws.Range("I16").Value = 4074.2
$ws.Range("K16").Value = 4074.2
$ws.Range("M16").Value = -3904.2
$ws.Range("H43").Value = 10000
$ws.Range("I43").Value = 10000
$ws.Range("K43").Value = 10000
$ws.Range("M43").Value = -9807
$ws.Range("H122").Value = 3641.2144
$ws.Range("I122").Value = 3600.2
$ws.Range("K122").Value = 10800.6
$ws.Range("M122").Value = -8350.599999999999
$ws.Range("H132").Value = 19446.143
$ws.Range("I132").Value = 19076.055
$ws.Range("J132").Value = 21666.666
$ws.Range("K132").Value = 57228.165
$ws.Range("L132").Value = 64999.99800000001
$ws.Range("M132").Value = -54698.165
$ws.Range("N132").Value = -70059.99800000001
$ws.Range("H136").Value = 3101.4
$ws.Range("I136").Value = 3101.4
$ws.Range("K136").Value = 9304.200000000001
$ws.Range("M136").Value = -6754.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 24333
$ws.Range("I39").Value = 24333
$ws.Range("K39").Value = 24333
$ws.Range("M39").Value = -23920
$ws.Range("H82").Value = 30300.5
$ws.Range("J82").Value = 30300.5
$ws.Range("L82").Value = 30300.5
$ws.Range("N82").Value = -31066.5
$ws.Range("H85").Value = 30300.5
$ws.Range("J85").Value = 30300.5
$ws.Range("L85").Value = 30300.5
$ws.Range("N85").Value = -32952.5
$ws.Range("H132").Value = 1383.1428
$ws.Range("I132").Value = 1383.1428
$ws.Range("K132").Value = 4149.428400000001
$ws.Range("M132").Value = -1619.428400000001
$ws.Range("H136").Value = 3252.4375
$ws.Range("I136").Value = 3135.9333
$ws.Range("K136").Value = 9407.7999
$ws.Range("M136").Value = -6857.7999
